$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.088.39'
$ws.Range("E2").Value = '  +5.19%  '
$ws.Range("D3").Value = '3.306.77'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'215.91"
$ws.Range("E5").Value = '  +2.13%  '
$ws.Range("D6").Value = "'633.02"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("D7").Value = "'0.423"
$ws.Range("E7").Value = '  +12.76%  '
$ws.Range("E8").Value = '  +8.12%  '
$ws.Range("D9").Value = "'0.997"
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("D10").Value = '3.297.08'
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("D11").Value = "'0.602"
$ws.Range("E11").Value = '  +3.96%  '
$ws.Range("D12").Value = "'0.0000274"
$ws.Range("E12").Value = '  +6.45%  '
$ws.Range("E13").Value = '  +1.89%  '
$ws.Range("D14").Value = "'34.65"
$ws.Range("E14").Value = '  +0.98%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.899.99'
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D16").Value = "'5.46"
$ws.Range("E16").Value = '  +3.23%  '
$ws.Range("D17").Value = '90.587.18'
$ws.Range("E17").Value = '  +4.81%  '
$ws.Range("D18").Value = '3.266.24'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").Value = "'3.24"
$ws.Range("E19").Value = '  +7.71%  '
$ws.Range("D20").Value = "'14.31"
$ws.Range("E20").Value = '  +1.34%  '
$ws.Range("D21").Value = "'434.68"
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("D22").Value = "'9.02"
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("B23").Value = 'Polkadot'
$ws.Range("C23").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D23").Value = "'5.39"
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").Value = "'0.0000187"
$ws.Range("E24").Value = '  +44.83%  '
$ws.Range("E25").Value = '  +6.54%  '
$ws.Range("D26").Value = "'12.22"
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").Value = '3.440.24'
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = "'76.33"
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").Value = "'0.178"
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("D32").Value = "'567.72"
$ws.Range("E32").Value = '  +3.67%  '
$ws.Range("D33").Value = "'8.72"
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("D34").Value = "'7.33"
$ws.Range("E34").Value = '  +5.10%  '
$ws.Range("D35").Value = "'1.38"
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").Value = "'1.93"
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = "'3.62"
$ws.Range("E37").Value = '  +22.92%  '
$ws.Range("D38").Value = "'22.85"
$ws.Range("E38").Value = '  +1.54%  '
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").Value = "'22.39"
$ws.Range("E40").Value = '  +3.68%  '
$ws.Range("D41").Value = "'0.995"
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").Value = "'0.397"
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("D43").Value = "'2.01"
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = "'186.04"
$ws.Range("E45").Value = '  +3.87%  '
$ws.Range("D46").Value = "'148.83"
$ws.Range("E46").Value = '  -5.64%  '
$ws.Range("D47").Value = "'44.35"
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("E48").Value = '  +9.01%  '
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").Value = "'0.633"
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = "'25.35"
$ws.Range("E51").Value = '  +4.63%  '

# Reset style on force-text cells back to the default/Normal style so
# we don't leave a stray "Text" number-format style behind (the source
# diff only touches cell values, not formatting).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
